$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Clear out the old content first so nothing stray is left behind.
# ---------------------------------------------------------------------------
$ws.Range("A1:D10").ClearContents()

# ---------------------------------------------------------------------------
# 2) Enter the date-looking text values using a leading apostrophe so they
#    are NOT auto-converted to real dates, then immediately reset the style
#    back to Normal so no quotePrefix/number-format residue survives once we
#    apply the real formatting in step 3.
# ---------------------------------------------------------------------------
$ws.Range("B3").Formula = "'2023-05-22"
$ws.Range("B3").Style = "Normal"

$ws.Range("B4").Formula = "'2023-05-22"
$ws.Range("B4").Style = "Normal"

$ws.Range("B6").Formula = "'2023-05-22"
$ws.Range("B6").Style = "Normal"

$ws.Range("B7").Formula = "'2023-05-22"
$ws.Range("B7").Style = "Normal"

# ---------------------------------------------------------------------------
# 3) Fill in the rest of the plain text values.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "kldfjvks"
$ws.Range("B2").Value = "tRIZTAN"
$ws.Range("C2").Value = "SDFSF"
$ws.Range("D2").Value = "SDFSFD"

$ws.Range("A3").Value = "tRIZTAN"
$ws.Range("C3").Value = "01:54:23"

$ws.Range("A4").Value = "tRIZTAN"
$ws.Range("C4").Value = "01:58:25"

$ws.Range("A5").Value = "tRIZTAN"

$ws.Range("A6").Value = "tRIZTAN"
$ws.Range("C6").Value = "02:03:57"

$ws.Range("A7").Value = "tRIZTAN"
$ws.Range("C7").Value = "02:04:25"

# ---------------------------------------------------------------------------
# 4) B5 holds a real datetime serial number with a custom display format.
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = 45068.08500893519
$ws.Range("B5").NumberFormat = "yyyy-mm-dd h:mm:ss"

# ---------------------------------------------------------------------------
# 5) Apply the (visually neutral) "general" horizontal alignment across the
#    A1:D4 block - this is what produces the shared style used by every
#    cell in that block.
# ---------------------------------------------------------------------------
$ws.Range("A1:D4").HorizontalAlignment = 1

# ---------------------------------------------------------------------------
# 6) Row heights for rows 1-4.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5
$ws.Rows.Item(4).RowHeight = 19.5

# ---------------------------------------------------------------------------
# 7) Column widths for columns A-D.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.29
$ws.Columns.Item(2).ColumnWidth = 13.29
$ws.Columns.Item(3).ColumnWidth = 13.29
$ws.Columns.Item(4).ColumnWidth = 13.29
